$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-22 Wednesday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-05-23 Thursday", 2) | Out-Null
$d.Content.Find.Execute("19+39=58", $true, $true, $false, $false, $false, $true, 1, $false, "11+32=43", 2) | Out-Null
$d.Content.Find.Execute("53+14=67", $true, $true, $false, $false, $false, $true, 1, $false, "47-28=19", 2) | Out-Null
$d.Content.Find.Execute("48-3=45", $true, $true, $false, $false, $false, $true, 1, $false, "4+50=54", 2) | Out-Null
$d.Content.Find.Execute("31+27=58", $true, $true, $false, $false, $false, $true, 1, $false, "98-91=7", 2) | Out-Null
$d.Content.Find.Execute("39+21=60", $true, $true, $false, $false, $false, $true, 1, $false, "49+7=56", 2) | Out-Null
$d.Content.Find.Execute("49+27=76", $true, $true, $false, $false, $false, $true, 1, $false, "96-4=92", 2) | Out-Null
$d.Content.Find.Execute("66-38=28", $true, $true, $false, $false, $false, $true, 1, $false, "0+91=91", 2) | Out-Null
$d.Content.Find.Execute("39-20=19", $true, $true, $false, $false, $false, $true, 1, $false, "7+56=63", 2) | Out-Null
$d.Content.Find.Execute("14+17=31", $true, $true, $false, $false, $false, $true, 1, $false, "20+59=79", 2) | Out-Null
$d.Content.Find.Execute("11+38=49", $true, $true, $false, $false, $false, $true, 1, $false, "64-16=48", 2) | Out-Null
$d.Content.Find.Execute("51-14=37", $true, $true, $false, $false, $false, $true, 1, $false, "44+30=74", 2) | Out-Null
$d.Content.Find.Execute("76-17=59", $true, $true, $false, $false, $false, $true, 1, $false, "55-18=37", 2) | Out-Null
$d.Content.Find.Execute("28+70=98", $true, $true, $false, $false, $false, $true, 1, $false, "2+81=83", 2) | Out-Null
$d.Content.Find.Execute("55-9=46", $true, $true, $false, $false, $false, $true, 1, $false, "13+48=61", 2) | Out-Null
$d.Content.Find.Execute("40+12=52", $true, $true, $false, $false, $false, $true, 1, $false, "65-56=9", 2) | Out-Null
$d.Content.Find.Execute("81+7=88", $true, $true, $false, $false, $false, $true, 1, $false, "90-7=83", 2) | Out-Null
$d.Content.Find.Execute("63+5=68", $true, $true, $false, $false, $false, $true, 1, $false, "16+19=35", 2) | Out-Null
$d.Content.Find.Execute("96-41=55", $true, $true, $false, $false, $false, $true, 1, $false, "24-14=10", 2) | Out-Null
$d.Content.Find.Execute("43+53=96", $true, $true, $false, $false, $false, $true, 1, $false, "34-9=25", 2) | Out-Null
$d.Content.Find.Execute("47-12=35", $true, $true, $false, $false, $false, $true, 1, $false, "75-61=14", 2) | Out-Null
$d.Content.Find.Execute("9+33=42", $true, $true, $false, $false, $false, $true, 1, $false, "81-0=81", 2) | Out-Null
$d.Content.Find.Execute("39+12=51", $true, $true, $false, $false, $false, $true, 1, $false, "49+30=79", 2) | Out-Null
$d.Content.Find.Execute("40+8=48", $true, $true, $false, $false, $false, $true, 1, $false, "50+1=51", 2) | Out-Null
$d.Content.Find.Execute("86-15=71", $true, $true, $false, $false, $false, $true, 1, $false, "2+9=11", 2) | Out-Null
$d.Content.Find.Execute("52+45=97", $true, $true, $false, $false, $false, $true, 1, $false, "8+32=40", 2) | Out-Null
$d.Content.Find.Execute("99-68=31", $true, $true, $false, $false, $false, $true, 1, $false, "83+13=96", 2) | Out-Null
$d.Content.Find.Execute("20+73=93", $true, $true, $false, $false, $false, $true, 1, $false, "18+34=52", 2) | Out-Null
$d.Content.Find.Execute("22+20=42", $true, $true, $false, $false, $false, $true, 1, $false, "82-56=26", 2) | Out-Null
$d.Content.Find.Execute("74-72=2", $true, $true, $false, $false, $false, $true, 1, $false, "72-54=18", 2) | Out-Null
$d.Content.Find.Execute("81-46=35", $true, $true, $false, $false, $false, $true, 1, $false, "60-49=11", 2) | Out-Null
$d.Content.Find.Execute("32+67=99", $true, $true, $false, $false, $false, $true, 1, $false, "24+52=76", 2) | Out-Null
$d.Content.Find.Execute("16+21=37", $true, $true, $false, $false, $false, $true, 1, $false, "72-33=39", 2) | Out-Null
$d.Content.Find.Execute("42-39=3", $true, $true, $false, $false, $false, $true, 1, $false, "85-59=26", 2) | Out-Null
$d.Content.Find.Execute("73+8=81", $true, $true, $false, $false, $false, $true, 1, $false, "96-52=44", 2) | Out-Null
$d.Content.Find.Execute("23+30=53", $true, $true, $false, $false, $false, $true, 1, $false, "23-9=14", 2) | Out-Null
$d.Content.Find.Execute("25+34=59", $true, $true, $false, $false, $false, $true, 1, $false, "10+67=77", 2) | Out-Null
$d.Content.Find.Execute("71-57=14", $true, $true, $false, $false, $false, $true, 1, $false, "73-23=50", 2) | Out-Null
$d.Content.Find.Execute("53-39=14", $true, $true, $false, $false, $false, $true, 1, $false, "37-19=18", 2) | Out-Null
$d.Content.Find.Execute("76-60=16", $true, $true, $false, $false, $false, $true, 1, $false, "74+24=98", 2) | Out-Null
$d.Content.Find.Execute("9+51=60", $true, $true, $false, $false, $false, $true, 1, $false, "55+24=79", 2) | Out-Null
$d.Content.Find.Execute("37+50=87", $true, $true, $false, $false, $false, $true, 1, $false, "3+48=51", 2) | Out-Null
$d.Content.Find.Execute("20-4=16", $true, $true, $false, $false, $false, $true, 1, $false, "19+37=56", 2) | Out-Null
$d.Content.Find.Execute("41-41=0", $true, $true, $false, $false, $false, $true, 1, $false, "80-3=77", 2) | Out-Null
$d.Content.Find.Execute("25+70=95", $true, $true, $false, $false, $false, $true, 1, $false, "94-9=85", 2) | Out-Null
$d.Content.Find.Execute("33-10=23", $true, $true, $false, $false, $false, $true, 1, $false, "51-37=14", 2) | Out-Null
$d.Content.Find.Execute("47+8=55", $true, $true, $false, $false, $false, $true, 1, $false, "91-61=30", 2) | Out-Null
$d.Content.Find.Execute("97-79=18", $true, $true, $false, $false, $false, $true, 1, $false, "7+16=23", 2) | Out-Null
$d.Content.Find.Execute("64-14=50", $true, $true, $false, $false, $false, $true, 1, $false, "25+41=66", 2) | Out-Null
$d.Content.Find.Execute("14+37=51", $true, $true, $false, $false, $false, $true, 1, $false, "18+6=24", 2) | Out-Null
$d.Content.Find.Execute("22-12=10", $true, $true, $false, $false, $false, $true, 1, $false, "31+17=48", 2) | Out-Null
$d.Content.Find.Execute("54-35=19", $true, $true, $false, $false, $false, $true, 1, $false, "24+38=62", 2) | Out-Null
$d.Content.Find.Execute("68+15=83", $true, $true, $false, $false, $false, $true, 1, $false, "59-46=13", 2) | Out-Null
$d.Content.Find.Execute("75-52=23", $true, $true, $false, $false, $false, $true, 1, $false, "39-19=20", 2) | Out-Null
$d.Content.Find.Execute("57+6=63", $true, $true, $false, $false, $false, $true, 1, $false, "1+74=75", 2) | Out-Null
$d.Content.Find.Execute("74-2=72", $true, $true, $false, $false, $false, $true, 1, $false, "62-36=26", 2) | Out-Null
$d.Content.Find.Execute("84+14=98", $true, $true, $false, $false, $false, $true, 1, $false, "57+1=58", 2) | Out-Null
$d.Content.Find.Execute("91-71=20", $true, $true, $false, $false, $false, $true, 1, $false, "30+48=78", 2) | Out-Null
$d.Content.Find.Execute("63+19=82", $true, $true, $false, $false, $false, $true, 1, $false, "71-6=65", 2) | Out-Null
$d.Content.Find.Execute("58-51=7", $true, $true, $false, $false, $false, $true, 1, $false, "0+15=15", 2) | Out-Null
$d.Content.Find.Execute("73-30=43", $true, $true, $false, $false, $false, $true, 1, $false, "14+68=82", 2) | Out-Null
$d.Content.Find.Execute("73-46=27", $true, $true, $false, $false, $false, $true, 1, $false, "88-4=84", 2) | Out-Null
$d.Content.Find.Execute("41+42=83", $true, $true, $false, $false, $false, $true, 1, $false, "92-32=60", 2) | Out-Null
$d.Content.Find.Execute("55-6=49", $true, $true, $false, $false, $false, $true, 1, $false, "53+28=81", 2) | Out-Null
$d.Content.Find.Execute("27+58=85", $true, $true, $false, $false, $false, $true, 1, $false, "93-10=83", 2) | Out-Null
$d.Content.Find.Execute("82-24=58", $true, $true, $false, $false, $false, $true, 1, $false, "73+19=92", 2) | Out-Null
$d.Content.Find.Execute("31+24=55", $true, $true, $false, $false, $false, $true, 1, $false, "6-2=4", 2) | Out-Null
$d.Content.Find.Execute("22+5=27", $true, $true, $false, $false, $false, $true, 1, $false, "33+19=52", 2) | Out-Null
$d.Content.Find.Execute("24-20=4", $true, $true, $false, $false, $false, $true, 1, $false, "86-4=82", 2) | Out-Null
$d.Content.Find.Execute("57-41=16", $true, $true, $false, $false, $false, $true, 1, $false, "88-35=53", 2) | Out-Null
$d.Content.Find.Execute("97-90=7", $true, $true, $false, $false, $false, $true, 1, $false, "94-3=91", 2) | Out-Null
$d.Content.Find.Execute("71+4=75", $true, $true, $false, $false, $false, $true, 1, $false, "36+6=42", 2) | Out-Null
$d.Content.Find.Execute("9+39=48", $true, $true, $false, $false, $false, $true, 1, $false, "62-18=44", 2) | Out-Null
$d.Content.Find.Execute("40+25=65", $true, $true, $false, $false, $false, $true, 1, $false, "26+7=33", 2) | Out-Null
$d.Content.Find.Execute("32+44=76", $true, $true, $false, $false, $false, $true, 1, $false, "29+14=43", 2) | Out-Null
$d.Content.Find.Execute("18+76=94", $true, $true, $false, $false, $false, $true, 1, $false, "26+16=42", 2) | Out-Null
$d.Content.Find.Execute("91-73=18", $true, $true, $false, $false, $false, $true, 1, $false, "4+53=57", 2) | Out-Null
$d.Content.Find.Execute("21+19=40", $true, $true, $false, $false, $false, $true, 1, $false, "20+24=44", 2) | Out-Null
$d.Content.Find.Execute("34+14=48", $true, $true, $false, $false, $false, $true, 1, $false, "32+49=81", 2) | Out-Null
$d.Content.Find.Execute("87-51=36", $true, $true, $false, $false, $false, $true, 1, $false, "79-0=79", 2) | Out-Null
$d.Content.Find.Execute("29+45=74", $true, $true, $false, $false, $false, $true, 1, $false, "84+6=90", 2) | Out-Null
$d.Content.Find.Execute("80-2=78", $true, $true, $false, $false, $false, $true, 1, $false, "45+4=49", 2) | Out-Null
$d.Content.Find.Execute("44+0=44", $true, $true, $false, $false, $false, $true, 1, $false, "13+11=24", 2) | Out-Null
$d.Content.Find.Execute("20+39=59", $true, $true, $false, $false, $false, $true, 1, $false, "18+4=22", 2) | Out-Null
$d.Content.Find.Execute("29-3=26", $true, $true, $false, $false, $false, $true, 1, $false, "92-7=85", 2) | Out-Null
$d.Content.Find.Execute("68-5=63", $true, $true, $false, $false, $false, $true, 1, $false, "48-7=41", 2) | Out-Null
$d.Content.Find.Execute("24-13=11", $true, $true, $false, $false, $false, $true, 1, $false, "36-31=5", 2) | Out-Null
$d.Content.Find.Execute("88-12=76", $true, $true, $false, $false, $false, $true, 1, $false, "90-36=54", 2) | Out-Null
$d.Content.Find.Execute("87-37=50", $true, $true, $false, $false, $false, $true, 1, $false, "27+21=48", 2) | Out-Null
$d.Content.Find.Execute("38+20=58", $true, $true, $false, $false, $false, $true, 1, $false, "11+5=16", 2) | Out-Null
$d.Content.Find.Execute("43-34=9", $true, $true, $false, $false, $false, $true, 1, $false, "91-6=85", 2) | Out-Null
$d.Content.Find.Execute("65-37=28", $true, $true, $false, $false, $false, $true, 1, $false, "6+59=65", 2) | Out-Null
$d.Content.Find.Execute("77-23=54", $true, $true, $false, $false, $false, $true, 1, $false, "34+45=79", 2) | Out-Null
$d.Content.Find.Execute("98-20=78", $true, $true, $false, $false, $false, $true, 1, $false, "76-70=6", 2) | Out-Null
$d.Content.Find.Execute("67-66=1", $true, $true, $false, $false, $false, $true, 1, $false, "85-6=79", 2) | Out-Null
$d.Content.Find.Execute("4+8=12", $true, $true, $false, $false, $false, $true, 1, $false, "4+79=83", 2) | Out-Null
$d.Content.Find.Execute("61-60=1", $true, $true, $false, $false, $false, $true, 1, $false, "80+17=97", 2) | Out-Null
$d.Content.Find.Execute("61-37=24", $true, $true, $false, $false, $false, $true, 1, $false, "17+18=35", 2) | Out-Null
$d.Content.Find.Execute("62-31=31", $true, $true, $false, $false, $false, $true, 1, $false, "53+35=88", 2) | Out-Null
$d.Content.Find.Execute("38+11=49", $true, $true, $false, $false, $false, $true, 1, $false, "33+52=85", 2) | Out-Null
$d.Content.Find.Execute("15+7=22", $true, $true, $false, $false, $false, $true, 1, $false, "71+18=89", 2) | Out-Null
